$wb = $excel.ActiveWorkbook

# "想去人数" (want-to-go count) increased on two rows, mirrored across
# the "展览" (sheet1) and "全部类型" (sheet4) worksheets, which carry
# duplicate data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 132
    $ws.Range("F3").Value = 82
}
